$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "Luca Perenzoni"
$ws.Range("B41").Value = "Riccardo Versini | Modium"
$ws.Range("C41").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("D41").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E41").Value = "Andrea Bertolini | Modium"
$ws.Range("F41").Value = "Alessio Koleci | FC. Stallions"
